# Update data: 5 August 2022
# Adds the newly published data point (date serial 44743 = 2022-07-01) to
# both the "Canada" sheet (row 32) and the "Province" sheet (rows 302-311,
# one row per province/territory) following the existing layout exactly.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Canada": append row 32
# ---------------------------------------------------------------------
$wsCanada = $wb.Worksheets.Item("Canada")

$lastRowCanada = 31
$newRowCanada = 32

$wsCanada.Range("A$newRowCanada").Value = 44743
$wsCanada.Range("A$newRowCanada").NumberFormat = $wsCanada.Range("A$lastRowCanada").NumberFormat

$wsCanada.Range("B$newRowCanada").Value = "Canada"
$wsCanada.Range("B$newRowCanada").NumberFormat = $wsCanada.Range("B$lastRowCanada").NumberFormat

$wsCanada.Range("D$newRowCanada").Value = 1007.1
$wsCanada.Range("E$newRowCanada").Value = 1165.3
$wsCanada.Range("C$newRowCanada").Formula = "=(D$newRowCanada-E$newRowCanada)/E$newRowCanada*100"

# ---------------------------------------------------------------------
# Sheet "Province": append rows 302-311, one per province/territory, in
# the same order used throughout the sheet.
# ---------------------------------------------------------------------
$wsProvince = $wb.Worksheets.Item("Province")

$provinceRows = @(
    @{ Row = 302; Name = "Newfoundland & Labrador"; D = 25.8;   E = 32.7;  FirstOfGroup = $true  },
    @{ Row = 303; Name = "Prince Edward Island";     D = 5.1;   E = 7;     FirstOfGroup = $false },
    @{ Row = 304; Name = "Nova Scotia";              D = 30;    E = 37.3;  FirstOfGroup = $false },
    @{ Row = 305; Name = "New Brunswick";            D = 28.3;  E = 34.1;  FirstOfGroup = $false },
    @{ Row = 306; Name = "Quebec";                   D = 185.1; E = 227.7; FirstOfGroup = $false },
    @{ Row = 307; Name = "Ontario";                  D = 426.9; E = 449.3; FirstOfGroup = $false },
    @{ Row = 308; Name = "Manitoba";                 D = 24.4;  E = 39.9;  FirstOfGroup = $false },
    @{ Row = 309; Name = "Saskatchewan";              D = 24;    E = 33.9;  FirstOfGroup = $false },
    @{ Row = 310; Name = "Alberta";                  D = 120.6; E = 178.2; FirstOfGroup = $false },
    @{ Row = 311; Name = "British Columbia";         D = 136.8; E = 125.3; FirstOfGroup = $false }
)

$templateStyledRow = 292   # first row of the previous date-group (has A & B styled)
$templatePlainRow  = 293   # following row of the previous date-group (only A styled)

foreach ($item in $provinceRows) {
    $r = $item.Row

    $wsProvince.Range("A$r").Value = 44743
    $wsProvince.Range("A$r").NumberFormat = $wsProvince.Range("A$templateStyledRow").NumberFormat

    $wsProvince.Range("B$r").Value = $item.Name
    if ($item.FirstOfGroup) {
        # Only the first row of a date-group carries the explicit "styled"
        # format (matches B132/B292/etc. in the existing data); the rest
        # stay on the sheet's default (unstyled) format, so NumberFormat is
        # intentionally left untouched for them.
        $wsProvince.Range("B$r").NumberFormat = $wsProvince.Range("B$templateStyledRow").NumberFormat
    }

    $wsProvince.Range("D$r").Value = $item.D
    $wsProvince.Range("E$r").Value = $item.E
    $wsProvince.Range("C$r").Formula = "=(D$r-E$r)/E$r*100"
}

# ---------------------------------------------------------------------
# Match the selection/active-cell state the source file ends up in
# (Canada -> A32, Province -> D312, Province tab active).
# ---------------------------------------------------------------------
$wsCanada.Range("A$newRowCanada").Select() | Out-Null
$wsCanada.Activate() | Out-Null

$wsProvince.Range("D312").Select() | Out-Null
$wsProvince.Activate() | Out-Null

Write-Output "Data updated: Canada row 32, Province rows 302-311"
